$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting existing rows 75-82 down to 76-83
$ws.Rows.Item(75).Insert()

# Columns A, B, D, E contain text that LOOKS numeric/date-like ("-548", "8/6/2025",
# "13", "ICD30326446 "). Force them to be stored as text (not auto-coerced to a
# number/date) by temporarily applying a text number format, then resetting the
# cell style back to Normal so no residual style index is left on the cell
# (matches the source workbook, where these cells carry no "s" attribute at all).
$textCols = @(1, 2, 4, 5)
foreach ($col in $textCols) {
    $ws.Cells.Item(75, $col).NumberFormat = "@"
}

$ws.Cells.Item(75, 1).Value = "-548"
$ws.Cells.Item(75, 2).Value = "8/6/2025"
$ws.Cells.Item(75, 3).Value = "Sucre 1533"
$ws.Cells.Item(75, 4).Value = "13"
$ws.Cells.Item(75, 5).Value = "ICD30326446 "
$ws.Cells.Item(75, 6).Value = "PEBCOM"
$ws.Cells.Item(75, 7).Value = "Pendiente"
$ws.Cells.Item(75, 8).Value = "Retirar columna"
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = "Desmonte"
$ws.Cells.Item(75, 11).Value = "Sin equipos"
$ws.Cells.Item(75, 12).Value = "Pasante"
$ws.Cells.Item(75, 13).Value = -58.44649
$ws.Cells.Item(75, 14).Value = -34.558808
$ws.Cells.Item(75, 15).Value = "Saavedra"
$ws.Cells.Item(75, 16).Value = "Capital Norte"

foreach ($col in $textCols) {
    $ws.Cells.Item(75, $col).Style = "Normal"
}
